# Fruta / hortaliza, semanal
# New weekly price observation was inserted before the existing row 147,
# pushing every subsequent record (old rows 147-206) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 147; Excel shifts rows 147:206 down to 148:207
# and the sheet dimension grows from A1:R206 to A1:R207 automatically.
$ws.Rows("147:147").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A147").Value = 10
$ws.Range("B147").Value = "Vega Modelo de Temuco"
$ws.Range("C147").Value = "La Araucanía"
$ws.Range("D147").Value = 44609
$ws.Range("E147").Value = 9
$ws.Range("F147").Value = 100112052
$ws.Range("G147").Value = "Albahaca"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 200
$ws.Range("K147").Value = 5000
$ws.Range("L147").Value = 5000
$ws.Range("M147").Value = 5000
$ws.Range("N147").Value = "$/paquete"
$ws.Range("O147").Value = "Región del Maule"
$ws.Range("P147").Value = 5000
$ws.Range("Q147").Value = 1
$ws.Range("R147").Value = "Hortaliza"
